$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")
$ws.Activate()

# Insert 3 new rows after RPA180_OfficeID1 (row 12) for OfficeID2/3/4
$ws.Rows("13:15").Insert()
$ws.Rows("13:15").RowHeight = 14.25

$ws.Range("A13").Value = "RPA180_OfficeID2"
$ws.Range("B13").Value = "RPA180_Amadeus_TJQ_TH_OfficeID2"

$ws.Range("A14").Value = "RPA180_OfficeID3"
$ws.Range("B14").Value = "RPA180_Amadeus_TJQ_TH_OfficeID3"

$ws.Range("A15").Value = "RPA180_OfficeID4"
$ws.Range("B15").Value = "RPA180_Amadeus_TJQ_TH_OfficeID4"

# Insert 1 new row after RPA183_OfficeID1 (now row 21) for OfficeID2
$ws.Rows("22:22").Insert()
$ws.Rows("22:22").RowHeight = 14.25

$ws.Range("A22").Value = "RPA183_OfficeID2"
$ws.Range("B22").Value = "RPA183_Amadeus_TJQ_MY_OfficeID2"

# Update selection to A10
$null = $ws.Range("A10").Select()
